# Work Order/System Setup: register the newest Engineering Item that was
# created in Salesforce ("Pro-PEItem-J7232") on the "Routing Master" sheet,
# which always tracks the most recently created Item Number / Salesforce Id
# pair for the downstream Work Order / Routing steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-J7232"
$ws.Range("D2").Value = "a345f000000uTj2AAE"
